$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Cells.Item(2, 6).Value = 30.32963275909424
$ws.Cells.Item(3, 6).Value = 29.87662672996521
$ws.Cells.Item(4, 6).Value = 30.04151892662048
$ws.Cells.Item(5, 6).Value = 29.95346665382385
$ws.Cells.Item(6, 6).Value = 29.8919312953949
$ws.Cells.Item(7, 6).Value = 29.99257564544678
$ws.Cells.Item(8, 6).Value = 29.8600332736969
$ws.Cells.Item(9, 6).Value = 29.96042943000793
$ws.Cells.Item(10, 6).Value = 29.90268731117249
$ws.Cells.Item(11, 6).Value = 30.15790510177612
$ws.Cells.Item(12, 6).Value = 29.94881272315979
$ws.Cells.Item(13, 6).Value = 29.84983134269714
$ws.Cells.Item(14, 6).Value = 29.84732055664062
$ws.Cells.Item(15, 6).Value = 29.85854697227478
$ws.Cells.Item(16, 6).Value = 29.84684300422668
$ws.Cells.Item(17, 6).Value = 29.90271091461182
$ws.Cells.Item(18, 6).Value = 29.78974080085754
$ws.Cells.Item(19, 6).Value = 29.87248206138611
$ws.Cells.Item(20, 6).Value = 29.97604894638061
$ws.Cells.Item(21, 6).Value = 30.10194325447083

$ws = $wb.Worksheets.Item("run_2")
$ws.Cells.Item(2, 6).Value = 30.09138011932373
$ws.Cells.Item(3, 6).Value = 30.02620077133179
$ws.Cells.Item(4, 6).Value = 29.78105735778809
$ws.Cells.Item(5, 6).Value = 29.9242033958435
$ws.Cells.Item(6, 6).Value = 29.88250732421875
$ws.Cells.Item(7, 6).Value = 30.01092576980591
$ws.Cells.Item(8, 6).Value = 29.96712827682495
$ws.Cells.Item(9, 6).Value = 29.9167549610138
$ws.Cells.Item(10, 6).Value = 29.87134194374084
$ws.Cells.Item(11, 6).Value = 30.20308613777161
$ws.Cells.Item(12, 6).Value = 29.82512235641479
$ws.Cells.Item(13, 6).Value = 29.88949131965637
$ws.Cells.Item(14, 6).Value = 29.84434723854065
$ws.Cells.Item(15, 6).Value = 29.95283985137939
$ws.Cells.Item(16, 6).Value = 29.81967234611511
$ws.Cells.Item(17, 6).Value = 29.9378604888916
$ws.Cells.Item(18, 6).Value = 29.86930465698243
$ws.Cells.Item(19, 6).Value = 29.85979056358337
$ws.Cells.Item(20, 6).Value = 29.8306450843811
$ws.Cells.Item(21, 6).Value = 30.18504309654236

$ws = $wb.Worksheets.Item("run_3")
$ws.Cells.Item(2, 6).Value = 30.11069369316101
$ws.Cells.Item(3, 6).Value = 29.97803521156311
$ws.Cells.Item(4, 6).Value = 29.77809429168701
$ws.Cells.Item(5, 6).Value = 29.88665294647217
$ws.Cells.Item(6, 6).Value = 29.91417837142944
$ws.Cells.Item(7, 6).Value = 30.06914973258972
$ws.Cells.Item(8, 6).Value = 29.80396938323974
$ws.Cells.Item(9, 6).Value = 29.90506839752197
$ws.Cells.Item(10, 6).Value = 29.93886971473694
$ws.Cells.Item(11, 6).Value = 30.18298029899597
$ws.Cells.Item(12, 6).Value = 29.83091020584106
$ws.Cells.Item(13, 6).Value = 29.9476249217987
$ws.Cells.Item(14, 6).Value = 29.84555625915528
$ws.Cells.Item(15, 6).Value = 29.93301701545716
$ws.Cells.Item(16, 6).Value = 29.86205339431763
$ws.Cells.Item(17, 6).Value = 29.85874438285828
$ws.Cells.Item(18, 6).Value = 29.92232370376587
$ws.Cells.Item(19, 6).Value = 29.91538381576538
$ws.Cells.Item(20, 6).Value = 29.83722257614136
$ws.Cells.Item(21, 6).Value = 30.12977623939514

$ws = $wb.Worksheets.Item("run_4")
$ws.Cells.Item(2, 6).Value = 30.14881658554077
$ws.Cells.Item(3, 6).Value = 29.90719413757324
$ws.Cells.Item(4, 6).Value = 29.84030437469482
$ws.Cells.Item(5, 6).Value = 30.01970434188843
$ws.Cells.Item(6, 6).Value = 29.85960102081299
$ws.Cells.Item(7, 6).Value = 29.99923586845398
$ws.Cells.Item(8, 6).Value = 29.84604287147522
$ws.Cells.Item(9, 6).Value = 29.9532253742218
$ws.Cells.Item(10, 6).Value = 29.86678409576416
$ws.Cells.Item(11, 6).Value = 30.1832709312439
$ws.Cells.Item(12, 6).Value = 29.86647367477417
$ws.Cells.Item(13, 6).Value = 29.98110914230347
$ws.Cells.Item(14, 6).Value = 29.82787919044494
$ws.Cells.Item(15, 6).Value = 29.91382050514221
$ws.Cells.Item(16, 6).Value = 29.86367154121399
$ws.Cells.Item(17, 6).Value = 29.97975444793701
$ws.Cells.Item(18, 6).Value = 29.81351137161255
$ws.Cells.Item(19, 6).Value = 29.92861342430115
$ws.Cells.Item(20, 6).Value = 30.05117678642273
$ws.Cells.Item(21, 6).Value = 30.09660768508911

$ws = $wb.Worksheets.Item("run_5")
$ws.Cells.Item(2, 6).Value = 30.05398321151733
$ws.Cells.Item(3, 6).Value = 29.93790721893311
$ws.Cells.Item(4, 6).Value = 29.76625204086304
$ws.Cells.Item(5, 6).Value = 29.84038472175598
$ws.Cells.Item(6, 6).Value = 29.87114834785461
$ws.Cells.Item(7, 6).Value = 30.12606835365296
$ws.Cells.Item(8, 6).Value = 30.06899333000183
$ws.Cells.Item(9, 6).Value = 30.03279876708984
$ws.Cells.Item(10, 6).Value = 29.90515494346619
$ws.Cells.Item(11, 6).Value = 30.23918104171753
$ws.Cells.Item(12, 6).Value = 29.83828186988831
$ws.Cells.Item(13, 6).Value = 29.93893051147461
$ws.Cells.Item(14, 6).Value = 29.87483358383179
$ws.Cells.Item(15, 6).Value = 30.01310229301453
$ws.Cells.Item(16, 6).Value = 29.8501181602478
$ws.Cells.Item(17, 6).Value = 29.94263219833374
$ws.Cells.Item(18, 6).Value = 29.89961314201355
$ws.Cells.Item(19, 6).Value = 29.87050485610962
$ws.Cells.Item(20, 6).Value = 29.82127785682678
$ws.Cells.Item(21, 6).Value = 30.24242830276489
